$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Keep gridlines / row-col headers visible, matching the authored sheet view
$excel.ActiveWindow.DisplayGridlines = $true
$excel.ActiveWindow.DisplayHeadings = $true

# ---------------------------------------------------------------------
# Drop the old 1-row header + 3 data rows; the report is being rebuilt
# with a title block on top and a wider 2010-2014 data table below.
# ---------------------------------------------------------------------
$ws.Range("A4:D7").Clear()

# ---------------------------------------------------------------------
# Row 2: report title "Balanza comercial" (merged A2:D2) - tan fill,
# thin black border framing the whole merged block.
# ---------------------------------------------------------------------
$ws.Range("A2:D2").Merge()
$ws.Range("A2:D2").Font.Size = 12
$ws.Range("A2:D2").Font.Bold = $false
$ws.Range("A2:D2").Font.Italic = $false
$ws.Range("A2:D2").Font.Strikethrough = $false
$ws.Range("A2:D2").Font.Underline = -4142
$ws.Range("A2:D2").Interior.Color = 13293535
$ws.Range("A2:D2").BorderAround(1, 2, -4142, 0)
$ws.Range("A2").Value = "Balanza comercial"

# ---------------------------------------------------------------------
# Row 3: "PAIS ORIGEN: ECUADOR" (merged A3:D3), small 8pt text
# ---------------------------------------------------------------------
$ws.Range("A3:D3").Merge()
$ws.Range("A3:D3").Font.Size = 8
$ws.Range("A3:D3").Font.Bold = $false
$ws.Range("A3:D3").Font.Italic = $false
$ws.Range("A3:D3").Font.Strikethrough = $false
$ws.Range("A3:D3").Font.Underline = -4142
$ws.Range("A3").Value = "PAIS ORIGEN: ECUADOR"

# ---------------------------------------------------------------------
# Row 4: "POSICION ARANCELARIA: ..." (merged A4:D4), small 8pt text
# ---------------------------------------------------------------------
$ws.Range("A4:D4").Merge()
$ws.Range("A4:D4").Font.Size = 8
$ws.Range("A4:D4").Font.Bold = $false
$ws.Range("A4:D4").Font.Italic = $false
$ws.Range("A4:D4").Font.Strikethrough = $false
$ws.Range("A4:D4").Font.Underline = -4142
$ws.Range("A4").Value = "POSICION ARANCELARIA: 09 CAFE TE YERBA MATE Y ESPECIAS."

# (row 5 intentionally stays blank, matching the authored layout)

# ---------------------------------------------------------------------
# Row 6: column headers - bold white text on dark-blue fill, centered,
# thin top border (same look the original single header row used).
# ---------------------------------------------------------------------
$ws.Range("A6").Value = "Periodo"
$ws.Range("B6").Value = "Total Impo CIF (US$)"
$ws.Range("C6").Value = "Total Expo FOB (US$)"
$ws.Range("D6").Value = "Balanza (US$)"
$ws.Range("A6:D6").Font.Bold = $true
$ws.Range("A6:D6").Font.Italic = $false
$ws.Range("A6:D6").Font.Strikethrough = $false
$ws.Range("A6:D6").Font.Underline = -4142
$ws.Range("A6:D6").Font.Color = 16777215
$ws.Range("A6:D6").Interior.Color = 8210719
$ws.Range("A6:D6").HorizontalAlignment = -4108
$ws.Range("A6:D6").Borders.Item(8).LineStyle = 1
$ws.Range("A6:D6").Borders.Item(8).Color = 0

# ---------------------------------------------------------------------
# Rows 7:11 - Impo/Expo/Balanza per year, 2010-2014
# ---------------------------------------------------------------------
$ws.Range("A7").Value = 2010
$ws.Range("B7").Value = 37795201.45
$ws.Range("C7").Value = 672517.58
$ws.Range("D7").Value = -37122683.87

$ws.Range("A8").Value = 2011
$ws.Range("B8").Value = 85960000.02
$ws.Range("C8").Value = 1107258.22
$ws.Range("D8").Value = -84852741.8

$ws.Range("A9").Value = 2012
$ws.Range("B9").Value = 63555667.7
$ws.Range("C9").Value = 1434370.56
$ws.Range("D9").Value = -62121297.14

$ws.Range("A10").Value = 2013
$ws.Range("B10").Value = 17327713.66
$ws.Range("C10").Value = 1282199.31
$ws.Range("D10").Value = -16045514.35

$ws.Range("A11").Value = 2014
$ws.Range("B11").Value = 9727020.47
$ws.Range("C11").Value = 815250.73
$ws.Range("D11").Value = -8911769.74

# ---------------------------------------------------------------------
# Final selection mirrors the authored workbook state
# ---------------------------------------------------------------------
$ws.Range("D6").Select()
